$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 804, shifting existing rows 804:873 down to 805:874.
$ws.Rows.Item(804).Insert()

# Populate the newly inserted row 804 with the new data record.
$ws.Cells.Item(804, 1).Value2 = 10
$ws.Cells.Item(804, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(804, 3).Value2 = "La Araucanía"
$ws.Cells.Item(804, 4).Value2 = 45132
$ws.Cells.Item(804, 5).Value2 = 9
$ws.Cells.Item(804, 6).Value2 = 100112045
$ws.Cells.Item(804, 7).Value2 = "Zapallo"
$ws.Cells.Item(804, 8).Value2 = "Camote"
$ws.Cells.Item(804, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(804, 10).Value2 = 850
$ws.Cells.Item(804, 11).Value2 = 500
$ws.Cells.Item(804, 12).Value2 = 500
$ws.Cells.Item(804, 13).Value2 = 500
$ws.Cells.Item(804, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(804, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(804, 16).Value2 = 500
$ws.Cells.Item(804, 17).Value2 = 1
$ws.Cells.Item(804, 18).Value2 = "Hortaliza"
